# Weekly price update: a new observation (row) was inserted at row 173
# of the data table, pushing the existing rows 173-202 down to 174-203.
# The new row 173 carries a new date (D) and volume (J) value; all the
# other columns repeat the constant "dimension" values used throughout
# this block of rows (same market / region / product / pricing unit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 173; Excel shifts rows 173:202 down to
# 174:203 and copies the surrounding row formatting (e.g. the date
# number format used in column D) onto the new row.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A173").Value = 3
$ws.Range("B173").Value = "Femacal de La Calera"
$ws.Range("C173").Value = "Coquimbo"
$ws.Range("D173").Value = 44505
$ws.Range("E173").Value = 5
$ws.Range("F173").Value = 100112039
$ws.Range("G173").Value = "Ciboulette"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 120
$ws.Range("K173").Value = 1500
$ws.Range("L173").Value = 1500
$ws.Range("M173").Value = 1500
$ws.Range("N173").Value = "$/docena de atados"
$ws.Range("O173").Value = "Provincia de Quillota"
$ws.Range("P173").Value = 500
$ws.Range("Q173").Value = 3
$ws.Range("R173").Value = "Hortaliza"
